$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 120.72
$ws.Range("I15").Value = 120.72
$ws.Range("K15").Value = 362.16
$ws.Range("M15").Value = -193.16

$ws.Range("H97").Value = 2500
$ws.Range("J97").Value = 3000
$ws.Range("L97").Value = 9000
$ws.Range("N97").Value = -9992

$ws.Range("H110").Value = 40276
$ws.Range("J110").Value = 40276
$ws.Range("L110").Value = 40276
$ws.Range("N110").Value = -48456

$ws.Range("H127").Value = 1597.1177
$ws.Range("I127").Value = 696.375
$ws.Range("J127").Value = 2397.7778
$ws.Range("K127").Value = 2089.125
$ws.Range("L127").Value = 7193.3334
$ws.Range("M127").Value = 2870.875
$ws.Range("N127").Value = -17113.3334

$ws.Range("H132").Value = 15213807
$ws.Range("I132").Value = 17929930
$ws.Range("J132").Value = 3513.5
$ws.Range("K132").Value = 53789790
$ws.Range("L132").Value = 10540.5
$ws.Range("M132").Value = -53787260
$ws.Range("N132").Value = -15600.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13308.684
$ws.Range("I32").Value = 9904.879000000001
$ws.Range("J32").Value = 17468.889
$ws.Range("K32").Value = 9904.879000000001
$ws.Range("L32").Value = 17468.889
$ws.Range("M32").Value = -9617.879000000001
$ws.Range("N32").Value = -18042.889

$ws.Range("H45").Value = 1858.6296
$ws.Range("I45").Value = 1022.3
$ws.Range("J45").Value = 2350.5881
$ws.Range("K45").Value = 1022.3
$ws.Range("L45").Value = 2350.5881
$ws.Range("M45").Value = -645.3
$ws.Range("N45").Value = -3104.5881

$ws.Range("H49").Value = 24999.5
$ws.Range("J49").Value = 24999.5
$ws.Range("L49").Value = 24999.5
$ws.Range("N49").Value = -25519.5

$ws.Range("H61").Value = 1528.7693
$ws.Range("I61").Value = 1038.9
$ws.Range("J61").Value = 3161.6667
$ws.Range("K61").Value = 1038.9
$ws.Range("L61").Value = 3161.6667
$ws.Range("M61").Value = -826.9000000000001
$ws.Range("N61").Value = -3585.6667

$ws.Range("H74").Value = 1684.1
$ws.Range("I74").Value = 1232.9166
$ws.Range("K74").Value = 1232.9166
$ws.Range("M74").Value = -358.9166

$ws.Range("H77").Value = 1684.1
$ws.Range("I77").Value = 1232.9166
$ws.Range("K77").Value = 6164.583000000001
$ws.Range("M77").Value = -1796.583000000001

$ws.Range("H122").Value = 2559.5557
$ws.Range("I122").Value = 1508.0714
$ws.Range("J122").Value = 6239.75
$ws.Range("K122").Value = 4524.2142
$ws.Range("L122").Value = 18719.25
$ws.Range("M122").Value = -2074.2142
$ws.Range("N122").Value = -23619.25

$ws.Range("H132").Value = 2120.1135
$ws.Range("I132").Value = 1024.6364
$ws.Range("K132").Value = 3073.9092
$ws.Range("M132").Value = -543.9092000000001

$ws.Range("H136").Value = 1528.7693
$ws.Range("I136").Value = 1038.9
$ws.Range("J136").Value = 3161.6667
$ws.Range("K136").Value = 3116.7
$ws.Range("L136").Value = 9485.000100000001
$ws.Range("M136").Value = -566.7000000000003
$ws.Range("N136").Value = -14585.0001

$ws.Range("H137").Value = 53275
$ws.Range("J137").Value = 53275
$ws.Range("L137").Value = 53275
$ws.Range("N137").Value = -63475

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 204.95833
$ws.Range("I80").Value = 109
$ws.Range("J80").Value = 273.5
$ws.Range("K80").Value = 109
$ws.Range("L80").Value = 273.5
$ws.Range("M80").Value = 889
$ws.Range("N80").Value = -2269.5

$ws.Range("H83").Value = 204.95833
$ws.Range("I83").Value = 109
$ws.Range("J83").Value = 273.5
$ws.Range("K83").Value = 545
$ws.Range("L83").Value = 1367.5
$ws.Range("M83").Value = 4447
$ws.Range("N83").Value = -11351.5

$ws.Range("H122").Value = 42069.332
$ws.Range("J122").Value = 42069.332
$ws.Range("L122").Value = 42069.332
$ws.Range("N122").Value = -51869.332

$ws.Range("H134").Value = 2855.5532
$ws.Range("I134").Value = 1511.2703
$ws.Range("K134").Value = 4533.810899999999
$ws.Range("M134").Value = -1998.810899999999

$ws.Range("H137").Value = 35418.57
$ws.Range("J137").Value = 35418.57
$ws.Range("L137").Value = 35418.57
$ws.Range("N137").Value = -45618.57

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3002.587
$ws.Range("I31").Value = 1378.8125
$ws.Range("J31").Value = 6714.0713
$ws.Range("K31").Value = 1378.8125
$ws.Range("L31").Value = 6714.0713
$ws.Range("M31").Value = -1083.8125
$ws.Range("N31").Value = -7304.0713

$ws.Range("H34").Value = 3002.587
$ws.Range("I34").Value = 1378.8125
$ws.Range("J34").Value = 6714.0713
$ws.Range("K34").Value = 1378.8125
$ws.Range("L34").Value = 6714.0713
$ws.Range("M34").Value = -1176.8125
$ws.Range("N34").Value = -7118.0713

$ws.Range("H68").Value = 70100.2
$ws.Range("J68").Value = 70100.2
$ws.Range("L68").Value = 70100.2
$ws.Range("N68").Value = -71598.2

$ws.Range("H71").Value = 70100.2
$ws.Range("J71").Value = 70100.2
$ws.Range("L71").Value = 210300.6
$ws.Range("N71").Value = -217788.6

$ws.Range("H99").Value = 5090.364
$ws.Range("I99").Value = 3078
$ws.Range("J99").Value = 6767.3335
$ws.Range("K99").Value = 3078
$ws.Range("L99").Value = 6767.3335
$ws.Range("M99").Value = -1580
$ws.Range("N99").Value = -9763.333500000001

$ws.Range("H126").Value = 5090.364
$ws.Range("I126").Value = 3078
$ws.Range("J126").Value = 6767.3335
$ws.Range("K126").Value = 9234
$ws.Range("L126").Value = 20302.0005
$ws.Range("M126").Value = -6764
$ws.Range("N126").Value = -25242.0005

$ws.Range("H132").Value = 5993.3335
$ws.Range("I132").Value = 6994.857
$ws.Range("J132").Value = 5117
$ws.Range("K132").Value = 20984.571
$ws.Range("L132").Value = 15351
$ws.Range("M132").Value = -18454.571
$ws.Range("N132").Value = -20411

$ws.Range("H134").Value = 6628.5454
$ws.Range("I134").Value = 7225.1763
$ws.Range("K134").Value = 21675.5289
$ws.Range("M134").Value = -19140.5289

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2415736
$ws.Range("J4").Value = 6826.5713
$ws.Range("L4").Value = 20479.7139
$ws.Range("N4").Value = -20703.7139

$ws.Range("H5").Value = 1901.7188
$ws.Range("I5").Value = 1043.3125
$ws.Range("J5").Value = 2760.125
$ws.Range("K5").Value = 3129.9375
$ws.Range("L5").Value = 8280.375
$ws.Range("M5").Value = -3017.9375
$ws.Range("N5").Value = -8504.375

$ws.Range("H6").Value = 478.83334
$ws.Range("I6").Value = 74.333336
$ws.Range("J6").Value = 883.3333
$ws.Range("K6").Value = 223.000008
$ws.Range("L6").Value = 2649.9999
$ws.Range("M6").Value = -110.000008
$ws.Range("N6").Value = -2875.9999

$ws.Range("H11").Value = 12736.077
$ws.Range("I11").Value = 12326.7
$ws.Range("J11").Value = 14100.667
$ws.Range("K11").Value = 36980.10000000001
$ws.Range("L11").Value = 42302.001
$ws.Range("M11").Value = -36840.10000000001
$ws.Range("N11").Value = -42582.001

$ws.Range("H64").Value = 4620.8
$ws.Range("I64").Value = 374.66666
$ws.Range("J64").Value = 10990
$ws.Range("K64").Value = 1123.99998
$ws.Range("L64").Value = 32970
$ws.Range("M64").Value = -853.9999800000001
$ws.Range("N64").Value = -33510

$ws.Range("H67").Value = 4620.8
$ws.Range("I67").Value = 374.66666
$ws.Range("J67").Value = 10990
$ws.Range("K67").Value = 1123.99998
$ws.Range("L67").Value = 32970
$ws.Range("M67").Value = -187.9999800000001
$ws.Range("N67").Value = -34842

$ws.Range("H109").Value = 6550
$ws.Range("I109").Value = 300
$ws.Range("J109").Value = 9675
$ws.Range("K109").Value = 900
$ws.Range("L109").Value = 29025
$ws.Range("M109").Value = 140
$ws.Range("N109").Value = -31105

$ws.Range("H113").Value = 531.5
$ws.Range("I113").Value = 544.4211
$ws.Range("J113").Value = 515.13336
$ws.Range("K113").Value = 1633.2633
$ws.Range("L113").Value = 1545.40008
$ws.Range("M113").Value = 536.7366999999999
$ws.Range("N113").Value = -5885.40008

$ws.Range("H131").Value = 8929545
$ws.Range("I131").Value = 27778808
$ws.Range("J131").Value = 947.1053000000001
$ws.Range("K131").Value = 83336424
$ws.Range("L131").Value = 2841.3159
$ws.Range("M131").Value = -83331384
$ws.Range("N131").Value = -12921.3159

$ws.Range("H135").Value = 1901.7188
$ws.Range("I135").Value = 1043.3125
$ws.Range("J135").Value = 2760.125
$ws.Range("K135").Value = 9389.8125
$ws.Range("L135").Value = 24841.125
$ws.Range("M135").Value = -6854.8125
$ws.Range("N135").Value = -29911.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2685.1853
$ws.Range("I102").Value = 2220
$ws.Range("K102").Value = 2220
$ws.Range("M102").Value = -598

$ws.Range("H122").Value = 3668.1667
$ws.Range("I122").Value = 2382.5
$ws.Range("J122").Value = 6239.5
$ws.Range("K122").Value = 7147.5
$ws.Range("L122").Value = 18718.5
$ws.Range("M122").Value = -4697.5
$ws.Range("N122").Value = -23618.5

$ws.Range("H132").Value = 2700.6667
$ws.Range("I132").Value = 1315.75
$ws.Range("K132").Value = 3947.25
$ws.Range("M132").Value = -1417.25

$ws.Range("H137").Value = 84752.5
$ws.Range("J137").Value = 84752.5
$ws.Range("L137").Value = 84752.5
$ws.Range("N137").Value = -94952.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2557.1428
$ws.Range("J46").Value = 2150
$ws.Range("L46").Value = 2150
$ws.Range("N46").Value = -2526

$ws.Range("H117").Value = 39800
$ws.Range("J117").Value = 39800
$ws.Range("L117").Value = 39800
$ws.Range("N117").Value = -48978

$ws.Range("H132").Value = 4283.423
$ws.Range("I132").Value = 1692.8572
$ws.Range("J132").Value = 7305.75
$ws.Range("K132").Value = 5078.571599999999
$ws.Range("L132").Value = 21917.25
$ws.Range("M132").Value = -2548.571599999999
$ws.Range("N132").Value = -26977.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6063904.5
$ws.Range("I132").Value = 3494.2942
$ws.Range("J132").Value = 15875997
$ws.Range("K132").Value = 10482.8826
$ws.Range("L132").Value = 47627991
$ws.Range("M132").Value = -7952.882599999999
$ws.Range("N132").Value = -47633051
